$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data in the desired final row order (rows 2-10), sorted descending by total_registros
$data = @(
    @("PAZ ANASTACIO JUANITA ROSA", 33),
    @("NIMA CARMEN KAREN DEL MILAGRO", 33),
    @("CARRION LAZARO MICHAEL LUIS", 30),
    @("ARRUNATEGUI ESPINOZA JOVANNY", 29),
    @("ALZAMORA CHERRES SIRLEY YASMIN", 28),
    @("ESPINOZA VALDIVIEZO JUNIOR RICARDO", 27),
    @("PULACHE LAZO VILMA YOHANA", 22),
    @("LILIAN ROXANA VEGA GARCÍA", 19),
    @("NAVARRO JUAREZ LIDIA", 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
